$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every timestamp in column A (rows 2-97) forward by exactly one day
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 + 1
}

# Update Notified Production (MW) values in column B for the affected rows
$bUpdates = @{
    18 = 7
    19 = 7
    20 = 7
    21 = 16
    22 = 36
    23 = 37
    24 = 59
    25 = 65
    26 = 211
    27 = 241
    28 = 281
    29 = 317
    30 = 699
    31 = 738
    32 = 833
    33 = 877
    34 = 1299
    35 = 1384
    36 = 1446
    37 = 1487
    38 = 1836
    39 = 1898
    40 = 1936
    41 = 1976
    42 = 2137
    43 = 2160
    44 = 2185
    45 = 2205
    46 = 2272
    47 = 2280
    48 = 2288
    49 = 2293
    50 = 2280
    51 = 2279
    52 = 2278
    53 = 2273
    54 = 2214
    55 = 2203
    56 = 2192
    57 = 2177
    58 = 2034
    59 = 2012
    60 = 1989
    61 = 1961
    62 = 1721
    63 = 1681
    64 = 1642
    65 = 1609
    66 = 1252
    67 = 1203
    68 = 1106
    69 = 1062
    70 = 617
    72 = 482
    73 = 447
    74 = 189
    75 = 137
    76 = 107
    77 = 94
    78 = 17
    79 = 16
    80 = 15
    81 = 15
    82 = 6
    83 = 6
    84 = 6
    85 = 6
}
foreach ($r in $bUpdates.Keys) {
    $ws.Cells.Item($r, 2).Value = $bUpdates[$r]
}
